# Apply "added guatemalan and swedish atlases" edit:
# - Clear the per-API "not working" markers for Guatemala (row 9) and
#   Sweden (row 12), leaving only the "later" (col D) and "not working"
#   (col F) markers, matching the plain bordered style used by already
#   empty cells in that row (style of G9/G12).
# - Remove the "Order of addition:" helper note (O15) and the GBIF /
#   Guatemala** / Sweden** helper cells (P15:P17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Guatemala row (9): clear C9, E9, H9, J9, K9, restoring the plain
#     bordered "empty" style (matching G9) instead of the "not working"
#     highlighted style.
$emptyStyleCell = $ws.Range("G9")
foreach ($addr in @("C9", "E9", "H9", "J9", "K9")) {
    $emptyStyleCell.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).ClearContents()
}

# --- Sweden row (12): same treatment as Guatemala, using G12 as the
#     style donor.
$emptyStyleCell12 = $ws.Range("G12")
foreach ($addr in @("C12", "E12", "H12", "J12", "K12")) {
    $emptyStyleCell12.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).ClearContents()
}

# --- Remove the "Order of addition:" note and its associated list.
$ws.Range("O15").ClearContents()
$ws.Range("P15").Clear()
$ws.Range("P16").Clear()
$ws.Range("P17").Clear()

# --- Restore the window to a non-minimized state and update the
#     visible top-left cell / selection to reflect where the edits were
#     made.
$ws.Application.ActiveWindow.WindowState = -4143
$ws.Range("A11").Select()
$ws.Range("O15:P17").Select()
